$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2189473333333334
$ws.Range("H2").Value = 0.656842
$ws.Range("I2").Value = 0.009402596261870986
$ws.Range("J2").Value = 0.009402596261870984
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 7.673425223617556
$ws.Range("R2").Value = 69.06082701255801
$ws.Range("S2").Value = 0.003596145075839217
$ws.Range("T2").Value = 0.003596145075839216
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2189473333333334
$ws.Range("H3").Value = 0.656842
$ws.Range("I3").Value = 0.009402596261870986
$ws.Range("J3").Value = 0.009402596261870984
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("Q3").Value = 6.549430478832667
$ws.Range("R3").Value = 58.944874309494
$ws.Range("S3").Value = 0.00306938576706449
$ws.Range("T3").Value = 0.003069385767064489
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2189473333333334
$ws.Range("H4").Value = 0.656842
$ws.Range("I4").Value = 0.009402596261870986
$ws.Range("J4").Value = 0.009402596261870984
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 1.951199008009556
$ws.Range("R4").Value = 17.560791072086
$ws.Range("S4").Value = 0.0009144279770967694
$ws.Range("T4").Value = 0.000914427977096769
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2189473333333334
$ws.Range("H5").Value = 0.656842
$ws.Range("I5").Value = 0.009402596261870986
$ws.Range("J5").Value = 0.009402596261870984
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 3.889129004812223
$ws.Range("R5").Value = 35.00216104331
$ws.Range("S5").Value = 0.00182263744187051
$ws.Range("T5").Value = 0.001822637441870509
$ws.Range("I6").Value = 0.8622887582286424
$ws.Range("J6").Value = 0.8622887582286423
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 703.7107755296613
$ws.Range("R6").Value = 6333.396979766952
$ws.Range("S6").Value = 0.3297935363267854
$ws.Range("T6").Value = 0.3297935363267853
$ws.Range("I7").Value = 0.8622887582286424
$ws.Range("J7").Value = 0.8622887582286423
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("S7").Value = 0.2814857479672377
$ws.Range("T7").Value = 0.2814857479672376
$ws.Range("I8").Value = 0.8622887582286424
$ws.Range("J8").Value = 0.8622887582286423
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 178.9396165499852
$ws.Range("R8").Value = 1610.456548949867
$ws.Range("S8").Value = 0.08385991942010727
$ws.Range("T8").Value = 0.08385991942010725
$ws.Range("I9").Value = 0.8622887582286424
$ws.Range("J9").Value = 0.8622887582286423
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 356.6623650267439
$ws.Range("R9").Value = 3209.961285240695
$ws.Range("S9").Value = 0.167149554514512
$ws.Range("T9").Value = 0.167149554514512
$ws.Range("G10").Value = 2.823530666666667
$ws.Range("H10").Value = 8.470592
$ws.Range("I10").Value = 0.1212552739852724
$ws.Range("J10").Value = 0.1212552739852723
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 98.95599598042311
$ws.Range("R10").Value = 890.603963823808
$ws.Range("S10").Value = 0.04637565458701341
$ws.Range("T10").Value = 0.04637565458701339
$ws.Range("G11").Value = 2.823530666666667
$ws.Range("H11").Value = 8.470592
$ws.Range("I11").Value = 0.1212552739852724
$ws.Range("J11").Value = 0.1212552739852723
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 84.46103236174933
$ws.Range("R11").Value = 760.1492912557439
$ws.Range("S11").Value = 0.03958260056971133
$ws.Range("T11").Value = 0.03958260056971132
$ws.Range("G12").Value = 2.823530666666667
$ws.Range("H12").Value = 8.470592
$ws.Range("I12").Value = 0.1212552739852724
$ws.Range("J12").Value = 0.1212552739852723
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 25.16253635981511
$ws.Range("R12").Value = 226.462827238336
$ws.Range("S12").Value = 0.01179240412058315
$ws.Range("T12").Value = 0.01179240412058315
$ws.Range("G13").Value = 2.823530666666667
$ws.Range("H13").Value = 8.470592
$ws.Range("I13").Value = 0.1212552739852724
$ws.Range("J13").Value = 0.1212552739852723
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 50.15395640828444
$ws.Range("R13").Value = 451.38560767456
$ws.Range("S13").Value = 0.02350461470796448
$ws.Range("T13").Value = 0.02350461470796447
$ws.Range("G14").Value = 0.1642436666666667
$ws.Range("H14").Value = 0.492731
$ws.Range("I14").Value = 0.007053371524214274
$ws.Range("J14").Value = 0.007053371524214274
$ws.Range("M14").Value = 35.04689966666667
$ws.Range("N14").Value = 105.140699
$ws.Range("O14").Value = 0.3824629895491901
$ws.Range("P14").Value = 0.3824629895491901
$ws.Range("Q14").Value = 5.756231306552111
$ws.Range("R14").Value = 51.80608175896901
$ws.Range("S14").Value = 0.002697653559552119
$ws.Range("T14").Value = 0.002697653559552119
$ws.Range("G15").Value = 0.1642436666666667
$ws.Range("H15").Value = 0.492731
$ws.Range("I15").Value = 0.007053371524214274
$ws.Range("J15").Value = 0.007053371524214274
$ws.Range("O15").Value = 0.3264402385872224
$ws.Range("P15").Value = 0.3264402385872223
$ws.Range("Q15").Value = 4.913064982546334
$ws.Range("R15").Value = 44.217584842917
$ws.Range("S15").Value = 0.002302504283208828
$ws.Range("T15").Value = 0.002302504283208827
$ws.Range("G16").Value = 0.1642436666666667
$ws.Range("H16").Value = 0.492731
$ws.Range("I16").Value = 0.007053371524214274
$ws.Range("J16").Value = 0.007053371524214274
$ws.Range("M16").Value = 8.911727666666666
$ws.Range("N16").Value = 26.735183
$ws.Range("O16").Value = 0.09725271102035077
$ws.Range("P16").Value = 0.09725271102035075
$ws.Range("Q16").Value = 1.463694828308111
$ws.Range("R16").Value = 13.173253454773
$ws.Range("S16").Value = 0.0006859595025635818
$ws.Range("T16").Value = 0.0006859595025635817
$ws.Range("G17").Value = 0.1642436666666667
$ws.Range("H17").Value = 0.492731
$ws.Range("I17").Value = 0.007053371524214274
$ws.Range("J17").Value = 0.007053371524214274
$ws.Range("M17").Value = 17.76285166666667
$ws.Range("N17").Value = 53.288555
$ws.Range("O17").Value = 0.1938440608432367
$ws.Range("P17").Value = 0.1938440608432367
$ws.Range("Q17").Value = 2.917435888189444
$ws.Range("R17").Value = 26.256922993705
$ws.Range("S17").Value = 0.001367254178889745
$ws.Range("T17").Value = 0.001367254178889745
